$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 507.8
$ws.Range("I94").Value = 484.75
$ws.Range("K94").Value = 484.75
$ws.Range("M94").Value = -33.75
$ws.Range("H132").Value = 4207.6416
$ws.Range("I132").Value = 4714.275
$ws.Range("J132").Value = 2648.7693
$ws.Range("K132").Value = 14142.825
$ws.Range("L132").Value = 7946.3079
$ws.Range("M132").Value = -11612.825
$ws.Range("N132").Value = -13006.3079
$ws.Range("H135").Value = 914.95
$ws.Range("I135").Value = 752.5789
$ws.Range("J135").Value = 4000
$ws.Range("K135").Value = 6773.2101
$ws.Range("L135").Value = 36000
$ws.Range("M135").Value = -4238.2101
$ws.Range("N135").Value = -41070
$ws.Range("H138").Value = 3045.7937
$ws.Range("I138").Value = 1703.425
$ws.Range("J138").Value = 5380.3477
$ws.Range("K138").Value = 5110.275
$ws.Range("L138").Value = 16141.0431
$ws.Range("M138").Value = 29.72500000000036
$ws.Range("N138").Value = -26421.0431

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1447.5
$ws.Range("I2").Value = 1070.7894
$ws.Range("K2").Value = 1070.7894
$ws.Range("M2").Value = -957.7893999999999
$ws.Range("H32").Value = 6129684
$ws.Range("I32").Value = 2825586.2
$ws.Range("K32").Value = 2825586.2
$ws.Range("M32").Value = -2825299.2
$ws.Range("H116").Value = 1447.5
$ws.Range("I116").Value = 1070.7894
$ws.Range("K116").Value = 1070.7894
$ws.Range("M116").Value = 1223.2106
$ws.Range("H118").Value = 99959
$ws.Range("J118").Value = 99959
$ws.Range("L118").Value = 99959
$ws.Range("N118").Value = -103273
$ws.Range("H140").Value = 149999
$ws.Range("J140").Value = 149999
$ws.Range("L140").Value = 149999
$ws.Range("N140").Value = -160359

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1447.5
$ws.Range("I3").Value = 1070.7894
$ws.Range("K3").Value = 1070.7894
$ws.Range("M3").Value = -956.7893999999999
$ws.Range("H64").Value = 570.8333
$ws.Range("I64").Value = 281.5
$ws.Range("J64").Value = 1149.5
$ws.Range("K64").Value = 281.5
$ws.Range("L64").Value = 1149.5
$ws.Range("M64").Value = -56.5
$ws.Range("N64").Value = -1599.5
$ws.Range("H67").Value = 570.8333
$ws.Range("I67").Value = 281.5
$ws.Range("J67").Value = 1149.5
$ws.Range("K67").Value = 281.5
$ws.Range("L67").Value = 1149.5
$ws.Range("M67").Value = 498.5
$ws.Range("N67").Value = -2709.5
$ws.Range("H102").Value = 26832.777
$ws.Range("I102").Value = 5091.8335
$ws.Range("J102").Value = 70314.664
$ws.Range("K102").Value = 5091.8335
$ws.Range("L102").Value = 70314.664
$ws.Range("M102").Value = -1846.8335
$ws.Range("N102").Value = -76804.664

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 162996
$ws.Range("J20").Value = 162996
$ws.Range("L20").Value = 162996
$ws.Range("N20").Value = -163468
$ws.Range("H30").Value = 162996
$ws.Range("J30").Value = 162996
$ws.Range("L30").Value = 162996
$ws.Range("N30").Value = -163178
$ws.Range("H128").Value = 162996
$ws.Range("J128").Value = 162996
$ws.Range("L128").Value = 162996
$ws.Range("N128").Value = -172956
$ws.Range("H132").Value = 2373.8696
$ws.Range("I132").Value = 2223.1428
$ws.Range("K132").Value = 6669.428400000001
$ws.Range("M132").Value = -4139.428400000001
$ws.Range("H134").Value = 2063.1052
$ws.Range("I134").Value = 1528.2188
$ws.Range("K134").Value = 4584.6564
$ws.Range("M134").Value = -2049.6564

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 7857.615
$ws.Range("I56").Value = 7857.615
$ws.Range("K56").Value = 7857.615
$ws.Range("M56").Value = -7327.615
$ws.Range("H68").Value = 1286.625
$ws.Range("I68").Value = 1258.6
$ws.Range("J68").Value = 1333.3334
$ws.Range("K68").Value = 3775.8
$ws.Range("L68").Value = 4000.0002
$ws.Range("M68").Value = -2964.8
$ws.Range("N68").Value = -5622.0002
$ws.Range("H71").Value = 1286.625
$ws.Range("I71").Value = 1258.6
$ws.Range("J71").Value = 1333.3334
$ws.Range("K71").Value = 11327.4
$ws.Range("L71").Value = 12000.0006
$ws.Range("M71").Value = -7271.4
$ws.Range("N71").Value = -20112.0006
$ws.Range("H113").Value = 1528.32
$ws.Range("I113").Value = 1636.8572
$ws.Range("J113").Value = 1486.1111
$ws.Range("K113").Value = 4910.571599999999
$ws.Range("L113").Value = 4458.3333
$ws.Range("M113").Value = -2740.571599999999
$ws.Range("N113").Value = -8798.3333
$ws.Range("H116").Value = 2007.25
$ws.Range("I116").Value = 2029
$ws.Range("J116").Value = 2000
$ws.Range("K116").Value = 6087
$ws.Range("L116").Value = 6000
$ws.Range("M116").Value = -2645
$ws.Range("N116").Value = -12884

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 5257500.5
$ws.Range("I7").Value = 5257500.5
$ws.Range("K7").Value = 5257500.5
$ws.Range("M7").Value = -5257388.5
$ws.Range("H8").Value = 5257500.5
$ws.Range("I8").Value = 5257500.5
$ws.Range("K8").Value = 5257500.5
$ws.Range("M8").Value = -5257361.5
$ws.Range("H11").Value = 21000000
$ws.Range("I11").Value = 23750000
$ws.Range("K11").Value = 23750000
$ws.Range("M11").Value = -23749861
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("H14").Value = 5750000
$ws.Range("I14").Value = 5750000
$ws.Range("K14").Value = 5750000
$ws.Range("M14").Value = -5749832
$ws.Range("H122").Value = 6500
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 6500
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 19500
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -24400
$ws.Range("H126").Value = 2332.5
$ws.Range("I126").Value = 1999.2858
$ws.Range("K126").Value = 5997.857400000001
$ws.Range("M126").Value = -3527.857400000001
$ws.Range("H134").Value = 71388
$ws.Range("J134").Value = 71388
$ws.Range("L134").Value = 214164
$ws.Range("N134").Value = -219234

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 24004
$ws.Range("I3").Value = 24004
$ws.Range("K3").Value = 24004
$ws.Range("M3").Value = -23892
$ws.Range("H7").Value = 2498.9167
$ws.Range("I7").Value = 2420.2
$ws.Range("K7").Value = 2420.2
$ws.Range("M7").Value = -2308.2
$ws.Range("H15").Value = 24004
$ws.Range("I15").Value = 24004
$ws.Range("K15").Value = 24004
$ws.Range("M15").Value = -23834
$ws.Range("H93").Value = 166668200
$ws.Range("I93").Value = 500000600
$ws.Range("K93").Value = 500000600
$ws.Range("M93").Value = -499999352
$ws.Range("H99").Value = 74894.836
$ws.Range("J99").Value = 106456.336
$ws.Range("L99").Value = 106456.336
$ws.Range("N99").Value = -112446.336
$ws.Range("H100").Value = 3599.25
$ws.Range("I100").Value = 2999
$ws.Range("K100").Value = 2999
$ws.Range("M100").Value = -2458
$ws.Range("H112").Value = 99894.5
$ws.Range("J112").Value = 99894.5
$ws.Range("L112").Value = 99894.5
$ws.Range("N112").Value = -102848.5
$ws.Range("H118").Value = 109930
$ws.Range("J118").Value = 109930
$ws.Range("L118").Value = 109930
$ws.Range("N118").Value = -113244
$ws.Range("H122").Value = 3865
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 3865
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 11595
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -16495
$ws.Range("H126").Value = 2498.9167
$ws.Range("I126").Value = 2420.2
$ws.Range("K126").Value = 7260.599999999999
$ws.Range("M126").Value = -4790.599999999999
$ws.Range("H132").Value = 3308.3582
$ws.Range("I132").Value = 2182.7646
$ws.Range("J132").Value = 3691.06
$ws.Range("K132").Value = 6548.293799999999
$ws.Range("L132").Value = 11073.18
$ws.Range("M132").Value = -4018.293799999999
$ws.Range("N132").Value = -16133.18

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 2348.0557
$ws.Range("I100").Value = 2258.923
$ws.Range("K100").Value = 4517.846
$ws.Range("M100").Value = -3976.846
$ws.Range("H116").Value = 99979
$ws.Range("J116").Value = 99979
$ws.Range("L116").Value = 99979
$ws.Range("N116").Value = -109157
$ws.Range("H122").Value = 2193.8125
$ws.Range("I122").Value = 2253.4
$ws.Range("K122").Value = 6760.200000000001
$ws.Range("M122").Value = -4310.200000000001
$ws.Range("H127").Value = 84940
$ws.Range("J127").Value = 84940
$ws.Range("L127").Value = 84940
$ws.Range("N127").Value = -94860
$ws.Range("H132").Value = 5491
$ws.Range("I132").Value = 4963.1304
$ws.Range("K132").Value = 14889.3912
$ws.Range("M132").Value = -12359.3912
